$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($c = 1; $c -le 19; $c++) {
    $t = $ws.Cells.Item(7, $c).Text
    Write-Host ("Col $c : " + $t)
}
